$d = $word.ActiveDocument
$r = $d.Content

function Replace-Text($search, $replace) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND: $search"
    }
}

# 1. Header contact line: add line break before LinkedIn, and append GitHub handle.
Replace-Text "princechopin@gmail.com | linkedin.com/in/chi-hung-wang-98334328" "princechopin@gmail.com | ^llinkedin.com/in/chi-hung-wang-98334328 | github.com/princechopin"

# 2. Add an extra space before "Core Competencies".
Replace-Text "history. Core" "history.  Core"

# 3. Core competencies list: trim "Software"/"Algorithms" qualifiers and merge the
#    "Cloud Computing Software Designs" paragraph into the preceding one.
Replace-Text "Design Automation Software| Computational Software | Business Intelligence | Data Visualization |^pCloud Computing Software Designs | Streaming | Big Data | Graph/AI/ML Algorithms " "Design Automation | Computational | Business Intelligence | Data Visualization | Cloud Computing | Streaming | Big Data | Graph/AI/ML"

# 4. Technical skills: annotate years of C/C++ experience.
Replace-Text "C/C++ | YACC" "C/C++(14) | YACC"

# 5. Work experience - Lingopal AI bullet 1: slash-separated tech list.
Replace-Text "CloudWatch client-server SQS" "CloudWatch/client-server/SQS"

# 6. Work experience - Lingopal AI bullet 2: "improve" -> "improving".
Replace-Text "translation quality, improve output" "translation quality, improving output"

# 7. Work experience - Lingopal AI bullet 3: slash-separated tool list + "Google translate".
Replace-Text "Interfacing with srt, ffmpeg, tsduck, Demucs, Whisper/Google|DeepL translator/ElevenLab APIs" "Interfacing with srt / ffmpeg / tsduck / Demucs / Whisper / Google translate / DeepL /ElevenLab APIs"

# 8. Work experience - Lingopal AI bullet 4: reorder "services"/"interfaces".
Replace-Text "Quickly built Python/CloudWatch interfaces for task-based GPU/CPU profiling/visualization services." "Quickly built Python/CloudWatch services for task-based GPU/CPU profiling/visualization."

# 9. Work experience - Lingopal AI bullet 5: "Refactor" -> "Refactored".
Replace-Text "Refactor automatic" "Refactored automatic"

# 10. Work experience - Lingopal AI bullet 6: add "real-time" and "translator".
Replace-Text "Building Tsduck streaming service for Lingopal." "Building Tsduck streaming service for real-time Lingopal translator."
